$d = $word.ActiveDocument

$replacements = @(
    @("923÷2=", "807÷8="),
    @("183÷4=", "171÷2="),
    @("127÷7=", "651÷2="),
    @("250÷9=", "701÷8="),
    @("442÷3=", "250÷8="),
    @("757÷4=", "759÷8="),
    @("303÷2=", "822÷3="),
    @("310÷6=", "267÷9="),
    @("878÷9=", "265÷9="),
    @("834÷2=", "116÷5="),
    @("244÷2=", "887÷5="),
    @("642÷2=", "371÷7="),
    @("895÷6=", "483÷4="),
    @("159÷7=", "969÷9="),
    @("482÷2=", "181÷5="),
    @("709÷2=", "938÷8="),
    @("702÷6=", "828÷4="),
    @("280÷3=", "539÷9="),
    @("852÷5=", "583÷9="),
    @("264÷2=", "195÷2="),
    @("310÷8=", "701÷8="),
    @("477÷6=", "114÷4="),
    @("479÷8=", "267÷6="),
    @("731÷7=", "209÷2="),
    @("519÷5=", "889÷2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
